$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Single-value cells (rows 1-12 in the summary block) ---
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "1208"
$t.Cell(6, 1).Range.Text  = "0.00260"
$t.Cell(7, 1).Range.Text  = "0.00016"
$t.Cell(8, 1).Range.Text  = "0.00009"
$t.Cell(9, 1).Range.Text  = "0.00024"
$t.Cell(10, 1).Range.Text = "0.00025"
$t.Cell(11, 1).Range.Text = "0.00032"
$t.Cell(12, 1).Range.Text = "0.19863"

# --- Rows that previously held a full tab-separated summary line, now
#     collapsed down to a single bare value ---
$t.Cell(44, 1).Range.Text = "99.93"
$t.Cell(45, 1).Range.Text = "0.2"
$t.Cell(46, 1).Range.Text = "278"

Write-Output "done"
